$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2:D9").Font.Color = 192
$ws.Range("F2:M9").VerticalAlignment = -4108
$ws.Range("F2:L5").MergeCells = $true
$ws.Range("F6:L7").MergeCells = $true
$ws.Range("F8:J8").MergeCells = $true
$ws.Range("F9:M9").MergeCells = $true
